$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text, preserving original style/number format
# (prevents Excel auto-converting numeric-looking strings like "1.000" or "23.394.30"
#  into actual numbers / dates, which would lose formatting and introduce float drift).
function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" "23.409.63"
$ws.Range("E2").Value = "  -0.23%  "
Set-TextValue "D3" "1.627.22"
Set-TextValue "D4" "0.9999"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("E5").Value = "  -0.17%  "
Set-TextValue "D6" "303.96"
$ws.Range("E6").Value = "  -1.48%  "
Set-TextValue "D7" "0.3787"
$ws.Range("E7").Value = "  +0.50%  "
Set-TextValue "D8" "52.12"
$ws.Range("E8").Value = "  -1.37%  "
Set-TextValue "D9" "0.3633"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D10" "0.08097"
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D11" "1.227"
$ws.Range("E11").Value = "  -3.73%  "
Set-TextValue "D12" "1.001"
$ws.Range("E12").Value = "  -0.17%  "
Set-TextValue "D13" "22.61"
$ws.Range("E13").Value = "  -2.24%  "
Set-TextValue "D14" "6.551"
$ws.Range("E14").Value = "  -1.52%  "
Set-TextValue "D15" "0.00001245"
$ws.Range("E15").Value = "  -2.83%  "
Set-TextValue "D16" "7.223"
$ws.Range("E16").Value = "  -3.15%  "
Set-TextValue "D17" "1.619.87"
$ws.Range("E17").Value = "  -0.93%  "
Set-TextValue "D18" "93.57"
$ws.Range("E18").Value = "  -1.17%  "
Set-TextValue "D19" "0.06904"
$ws.Range("E19").Value = "  -0.65%  "
Set-TextValue "D20" "17.90"
$ws.Range("E20").Value = "  -2.53%  "
Set-TextValue "D21" "0.9998"
$ws.Range("E21").Value = "  -0.30%  "
Set-TextValue "D22" "6.409"
$ws.Range("E22").Value = "  -2.51%  "
Set-TextValue "D23" "23.412.84"
$ws.Range("E23").Value = "  -0.26%  "
Set-TextValue "D24" "12.72"
$ws.Range("E24").Value = "  -1.37%  "
Set-TextValue "D25" "3.224"
$ws.Range("E25").Value = "  +3.42%  "
Set-TextValue "D26" "2.436"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("E28").Value = "  -1.84%  "
Set-TextValue "D29" "5.293"
$ws.Range("E29").Value = "  -0.62%  "
Set-TextValue "D30" "134.38"
$ws.Range("E30").Value = "  -1.27%  "
$ws.Range("E31").Value = "  -4.76%  "
Set-TextValue "D32" "1.802.68"
$ws.Range("E32").Value = "  -0.77%  "
Set-TextValue "D33" "6.797"
$ws.Range("E33").Value = "  -0.08%  "
Set-TextValue "D34" "11.03"
$ws.Range("E34").Value = "  +5.52%  "
Set-TextValue "D35" "0.9524"
$ws.Range("E35").Value = "  -2.34%  "
Set-TextValue "D36" "0.02796"
$ws.Range("E36").Value = "  -0.61%  "
Set-TextValue "D37" "0.2532"
$ws.Range("E37").Value = "  -0.27%  "
Set-TextValue "D38" "0.08816"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D39" "0.07199"
$ws.Range("E39").Value = "  -3.13%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D40" "6.104"
$ws.Range("E40").Value = "  -1.95%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D41" "0.7073"
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D42" "1.356"
$ws.Range("E42").Value = "  -2.86%  "
Set-TextValue "D43" "16.17"
$ws.Range("E43").Value = "  +0.10%  "
Set-TextValue "D44" "12.33"
$ws.Range("E44").Value = "  -2.48%  "
Set-TextValue "D45" "0.6469"
$ws.Range("E45").Value = "  -2.27%  "
Set-TextValue "D46" "2.327"
$ws.Range("E46").Value = "  -1.29%  "
Set-TextValue "D47" "0.9996"
$ws.Range("E47").Value = "  -0.25%  "
Set-TextValue "D48" "3.994"
$ws.Range("E48").Value = "  -1.23%  "
Set-TextValue "D49" "0.07991"
$ws.Range("E49").Value = "  -0.57%  "
Set-TextValue "D50" "1.204"
Set-TextValue "D51" "125.67"
$ws.Range("E51").Value = "  -4.00%  "
